$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62, shifting existing rows 62-66 down to 63-67
$ws.Rows.Item(62).Insert()

# Copy the date cell's number format (only, not the whole row) from the row below
# (previously row 62, now shifted to row 63) into the new D62 cell
$ws.Range("D63").Copy()
$ws.Range("D62").PasteSpecial(-4122) # xlPasteFormats

# Fill in the constant columns for the new row (same as adjacent rows)
$ws.Range("A62").Value = 2
$ws.Range("B62").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C62").Value = "Coquimbo"
$ws.Range("D62").Value = 44783
$ws.Range("E62").Value = 4
$ws.Range("F62").Value = 100112022
$ws.Range("G62").Value = "Arveja Verde"
$ws.Range("H62").Value = "Perfection"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 400
$ws.Range("K62").Value = 27000
$ws.Range("L62").Value = 29000
$ws.Range("M62").Value = 28000
$ws.Range("N62").Value = "$/malla 25 kilos"
$ws.Range("O62").Value = "Provincia de Limarí"
$ws.Range("P62").Value = 1120
$ws.Range("Q62").Value = 25
$ws.Range("R62").Value = "Hortaliza"
